# Trackers.xlsx - Add LinkItem events to trackers.js (CMS-16008)
# Adds descriptive "comment" cells in column B for a few existing event
# sections, and adds a brand new "cms_linkItem" tracker table (rows 34-39)
# on the "CMS" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1) New description cells (column B) for existing sections
# ---------------------------------------------------------------------
$ws.Range("B23").Value = "When edit mode loads, changes view or heartbeat every minute"
$ws.Range("B27").Value = "When content is saved while editing"
$ws.Range("B30").Value = "When a user clicks something in the UI, the action specify what and where."

# ---------------------------------------------------------------------
# 2) New "cms_linkItem" table, rows 34-39
# ---------------------------------------------------------------------

# Row 34: section title + description, same layout/style as row 12
# (A: s3 bold-ish header style, B: s2 plain text)
$ws.Range("A12:B12").Copy($ws.Range("A34"))
$ws.Range("A34").Value = "cms_linkItem"
$ws.Range("B34").Value = "LinkItem; Feature analysis on how and if LinkItem addOn is used."

# Row 35: property row, same layout/style as row 15
# (A: s4, B: s2, C: s5, D: s2)
$ws.Range("A15:D15").Copy($ws.Range("A35"))
$ws.Range("A35").Value = "commandType"
$ws.Range("B35").Value = "String"
$ws.Range("C35").Value = """init"" "
$ws.Range("D35").Value = "Initiated if the client has LinkItem add-on"

# Rows 36-39: value rows, same layout/style as row 16 (C: s5, D: s2)
$ws.Range("C16:D16").Copy($ws.Range("C36"))
$ws.Range("C36").Value = """drop"""
$ws.Range("D36").Value = "When user drops a link into the LinkItem"

$ws.Range("C16:D16").Copy($ws.Range("C37"))
$ws.Range("C37").Value = """click"""
$ws.Range("D37").Value = "When user clicks on link item create link"

$ws.Range("C16:D16").Copy($ws.Range("C38"))
$ws.Range("C38").Value = """dialogSave"""
$ws.Range("D38").Value = "When user click save in dialog"

$ws.Range("C16:D16").Copy($ws.Range("C39"))
$ws.Range("C39").Value = """clear"""
$ws.Range("D39").Value = "When user clears the LinkItem"

# New fill/style introduced for the "cms_linkItem" property values
# (light grey FFE7E6E6, matching the new cellXfs entry added upstream)
$ws.Range("C35:C39").Interior.Color = 15132391

# ---------------------------------------------------------------------
# 3) Update sheet view / selection to match the edited document
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F43").Select()
